# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 6b0570f6... row on both the zh-cn and de-de report
# sheets, reflecting the regenerated handback report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D6").Value = "2016-03-09 04:25:47"
$zhcn.Range("D7").Value = "2016-03-09 04:25:47"
$zhcn.Range("G6").Value = "2016-03-09 04:26:43"
$zhcn.Range("G7").Value = "2016-03-09 04:26:43"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D6").Value = "2016-03-09 04:25:50"
$dede.Range("D7").Value = "2016-03-09 04:25:50"
$dede.Range("G6").Value = "2016-03-09 04:26:48"
$dede.Range("G7").Value = "2016-03-09 04:26:48"
